# Dash DI - V1
# Keep only the "Debêntures" data row (originally row 4) and remove the
# other "Títulos Privados" data rows (originally rows 2, 3, 5, 6, 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 7 (the 2nd..6th data rows), leaving only the header
# row (row 1) and the first data row (row 2).
$ws.Range("A3:G7").EntireRow.Delete() | Out-Null

# Overwrite what remains of the single data row (row 2) with the values
# that used to live in row 4 ("Debêntures" / HAPVIDA / JERA2026 / ...).
$ws.Range("A2").Value = "Debêntures"
$ws.Range("B2").Value = "HAPVIDA PARTICIPACOES E INVESTIMENTOS S/A"
$ws.Range("C2").Value = "JERA2026"
$ws.Range("D2").Value = "HAPV21"
$ws.Range("E2").Value = "BRHAPVDBS014"
$ws.Range("F2").Value = "362409-HAP"
$ws.Range("G2").Value = "2026-07-10 00:00:00"
